{"js": "// The edit: a single character \"d\" was typed at the very start of the\n// document (before the existing \"Module- ALTP\" text). Word's automatic\n// \"_GoBack\" bookmark (which marks the location of the user's last edit)\n// therefore moves from the end of the \"Flie c\u00e2u h\u1ecfi l\u00e0 k \u0111c x\u00f3a\" paragraph\n// to right after the newly typed \"d\" in the first paragraph.\n\n// 1) Remove the existing \"_GoBack\" bookmark (it currently sits at the end\n//    of the \"Flie c\u00e2u h\u1ecfi l\u00e0 k \u0111c x\u00f3a\" paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Type \"d\" at the very beginning of the first paragraph (\"Module- ALTP\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst insertedRange = firstParagraph.getRange(\"Start\").insertText(\"d\", Word.InsertLocation.before);\nawait context.sync();\n\n// 3) Re-create \"_GoBack\" as a collapsed bookmark immediately after the \"d\"\n//    that was just typed, matching Word's behaviour of tracking the last\n//    edit position.\nconst caretAfterD = insertedRange.getRange(\"End\");\ncaretAfterD.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The edit: a single character \"d\" was typed at the very start of the\n# document (before the existing \"Module- ALTP\" text). Word's automatic\n# \"_GoBack\" bookmark (which marks the location of the user's last edit)\n# therefore moves from the end of the \"Flie c\u00e2u h\u1ecfi l\u00e0 k \u0111c x\u00f3a\" paragraph\n# to right after the newly typed \"d\" in the first paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (currently at the end of the\n#    \"Flie c\u00e2u h\u1ecfi l\u00e0 k \u0111c x\u00f3a\" paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Type \"d\" at the very beginning of the first paragraph (\"Module- ALTP\").\n$firstParagraph = $d.Paragraphs.Item(1).Range\n$firstParagraph.InsertBefore(\"d\")\n\n# 3) Re-create \"_GoBack\" as a collapsed bookmark immediately after the \"d\"\n#    that was just typed, matching Word's behaviour of tracking the last\n#    edit position.\n$caretAfterD = $d.Range(0, 1)\n$caretAfterD.Collapse(0)   # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $caretAfterD)\n"}
